# Change the year in the astromap link from 2018 to 2022, and collapse
# the paragraph's three differently-formatted runs (plain prefix text +
# hyperlink-styled URL + plain suffix) into a single, unformatted run
# (preceded by an empty run), matching the author's manual edit.

$d = $word.ActiveDocument

# Locate the target paragraph ("Mapky v tomto dokumente pripravil Jan
# Hollan, CzechGlobe (...).") via Find, then expand to the whole paragraph.
$rng = $d.Content
$found = $rng.Find.Execute(
    "Mapky v tomto dokumente pripravil Jan Hollan", $false, $false, $false,
    $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Mapky v tomto dokumente ...' paragraph"
}
$rng.Expand(4) | Out-Null  # wdParagraph -> whole paragraph, including the mark

# Work on a copy that excludes the trailing paragraph mark.
$pr = $d.Range($rng.Start, $rng.End - 1)

$newText = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe " + `
    "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

# Reduce the existing (3-run, variously-formatted) content down to a single
# placeholder character first. Doing this keeps the range non-empty so the
# following InsertXML splices new runs in place of it, instead of treating
# a zero-length range as "the whole (now-empty) paragraph" and clobbering
# the paragraph's own properties (borders/paraId/etc.).
$pr.Text = "X"
$target = $d.Range($pr.Start, $pr.Start + 1)

# Insert the replacement as raw OOXML so the result is an empty run
# followed by one plain run holding the full sentence, with no leftover
# character formatting (rFonts/sz/lang/Hyperlink style) on it.
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body>' + `
    '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xmlFrag) | Out-Null
